# Feedback from industry visits week#7
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new feedback text for row 27 (B27). Note line 2 contains a
# right single quotation mark (U+2019) in "don't", unlike the straight
# apostrophes used everywhere else in the passage.
$rsquo = [char]0x2019
$lines = @(
    "1.plan ahead for what you don't know yet, keep an eye on elements you forsee could become an issue down the line (Graham & Ricardo - Outplay)",
    "2.Developers often create 'rods for their own back,' making things harder for themselves when they don${rsquo}t need to, leaving things in they are attached to too long for example (Graham & Ricardo - Outplay) ",
    "3.In these chats we have as little as 5 minutes to garner as much advice and ideas as we can, use their experience and time wisely (Ken)",
    "4.Avoid perfectionist behaviour, sometimes you need to get something to a stage that is 'good enough for now,' and improve upon it later if time allows, otherwise the game won't become finished. 'Still continues to work - Doesn't break' (Graham & Ricardo - Outplay)",
    "5.Proof of concept as quickly as possible, more time to improve upon it! (Graham & Ricardo - Outplay) "
)
$feedbackText = ($lines -join "`n") + "`n"

$cell = $ws.Range("B27")
$cell.Value = $feedbackText
$cell.WrapText = $true
$cell.Borders.LineStyle = 1

# Give row 27 the same 210pt height it needs to show the new note.
$ws.Rows.Item(27).RowHeight = 210

# Move the active selection/view the way the author left it.
$ws.Range("B28").Select()
